$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend the Excel Table (ListObject) "Tableau1" by one empty row (A1:G7 -> A1:G8)
$table = $ws.ListObjects.Item("Tableau1")
$table.ListRows.Add() | Out-Null

# New content below the table
$ws.Range("A9").Value = "Branche multithread"

$ws.Range("A12").Value = "FAIRE UNE COURBE!!!"
$ws.Range("A12").Interior.ThemeColor = 4
$ws.Range("A12").Interior.TintAndShade = 0.79998168889431442
$ws.Range("A12").Borders.Item(7).LineStyle = 1
$ws.Range("A12").Borders.Item(7).Color.ThemeColor = 4
$ws.Range("A12").Borders.Item(7).TintAndShade = 0.39997558519241921
$ws.Range("A12").Borders.Item(8).LineStyle = 1
$ws.Range("A12").Borders.Item(8).Color.ThemeColor = 4
$ws.Range("A12").Borders.Item(8).TintAndShade = 0.39997558519241921
$ws.Range("A12").Borders.Item(9).LineStyle = 1
$ws.Range("A12").Borders.Item(9).Color.ThemeColor = 4
$ws.Range("A12").Borders.Item(9).TintAndShade = 0.39997558519241921

$ws.Range("A12").Select()
